# Applies the OOXML diff:
#  - Student sheet: appends a trailing "," inside the CONCATENATE() formulas
#    in F2:F4, and re-selects F2:F4 / makes Student the active (selected) tab.
#  - DegreeCredit sheet: appends a trailing "," inside the "&"-concatenation
#    formulas in D2:D49, and re-selects D2:D49.
#  - DegreePlan sheet: is no longer the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- Student sheet -----------------------------------------------------
$student = $wb.Worksheets.Item("Student")

$studentFormula = '=CONCATENATE("new Student{StudentID = ",A:A," , Family = ''''",B:B,"'''' , Given =''''",C:C,"''''",", Snumber =",D:D," , number_919=",E:E,"},")'

$student.Range("F2").Formula = $studentFormula
$student.Range("F3").Formula = $studentFormula
$student.Range("F4").Formula = $studentFormula

# --- DegreeCredit sheet --------------------------------------------------
$degreeCredit = $wb.Worksheets.Item("DegreeCredit")

for ($row = 2; $row -le 49; $row++) {
    $f = '=" new DegreeCredit{"&$A$1&" = "&A' + $row + '&" , "&$B$1&" = "&B' + $row + '&" ,"&$C$1&" ="&C' + $row + '&"},"'
    $degreeCredit.Range("D$row").Formula = $f
}

# --- Selections / active tab ---------------------------------------------
# Make DegreeCredit's selection D2:D49 (active cell D2)
$degreeCredit.Range("D2:D49").Select() | Out-Null

# Make Student the active sheet with F2:F4 selected (active cell F2) and
# tabSelected -- this also shifts the workbook's activeTab away from
# DegreePlan (which previously had tabSelected="1").
$student.Activate() | Out-Null
$student.Range("F2:F4").Select() | Out-Null
